$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 447, shifting the existing data
# (old rows 447-513) down to rows 448-514.
$ws.Rows.Item(447).Insert()

# Populate the newly inserted row 447 with the new record.
$ws.Range("A447").Value = 10
$ws.Range("B447").Value = "Vega Modelo de Temuco"
$ws.Range("C447").Value = "La Araucanía"
$ws.Range("D447").Value = 44522
$ws.Range("E447").Value = 9
$ws.Range("F447").Value = 100112021
$ws.Range("G447").Value = "Ají"
$ws.Range("H447").Value = "Americana (o)"
$ws.Range("I447").Value = "Primera"
$ws.Range("J447").Value = 50
$ws.Range("K447").Value = 30000
$ws.Range("L447").Value = 30000
$ws.Range("M447").Value = 30000
$ws.Range("N447").Value = "$/caja 25 kilos"
$ws.Range("O447").Value = "Provincia de Limarí"
$ws.Range("P447").Value = 1200
$ws.Range("Q447").Value = 25
$ws.Range("R447").Value = "Hortaliza"
